$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (prices) are written as text, not auto-converted numbers,
# matching the source workbook which stores all these cells as inline strings.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.112.16'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.488.85'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.81'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.99'
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.525'
$ws.Range('E7').Value = '  +2.22%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.13'
$ws.Range('E10').Value = '  +7.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.48'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.19'
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.879.64'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.494.99'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.855'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.063.00'
$ws.Range('E18').Value = '  +4.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.94'
$ws.Range('E19').Value = '  +5.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('E20').Value = '  +5.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0939'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.68'
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.46'
$ws.Range('E23').Value = '  +8.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '250.35'
$ws.Range('E25').Value = '  +3.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.19'
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  +4.04%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('E29').Value = '  -3.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.91'
$ws.Range('E30').Value = '  +6.36%  '
$ws.Range('E31').Value = '  +9.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.40'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.67'
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.47'
$ws.Range('E34').Value = '  +4.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0791'
$ws.Range('E35').Value = '  +4.23%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.98'
$ws.Range('E37').Value = '  +6.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('E38').Value = '  +6.14%  '
$ws.Range('E39').Value = '  +3.74%  '
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.43'
$ws.Range('E41').Value = '  -4.73%  '
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.12'
$ws.Range('E43').Value = '  +2.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0298'
$ws.Range('E44').Value = '  +3.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.968.06'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.00'
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.05'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.26'
$ws.Range('E50').Value = '  +9.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.39'
$ws.Range('E51').Value = '  +3.34%  '
